$d = $word.ActiveDocument

$replacements = @(
    @{old = "725÷5="; new = "663÷9="},
    @{old = "510÷5="; new = "283÷4="},
    @{old = "610÷4="; new = "933÷2="},
    @{old = "498÷9="; new = "703÷7="},
    @{old = "261÷2="; new = "964÷8="},
    @{old = "198÷4="; new = "918÷7="},
    @{old = "165÷2="; new = "707÷2="},
    @{old = "979÷6="; new = "261÷9="},
    @{old = "230÷6="; new = "722÷9="},
    @{old = "703÷4="; new = "935÷5="},
    @{old = "882÷6="; new = "969÷8="},
    @{old = "997÷3="; new = "125÷3="},
    @{old = "146÷4="; new = "556÷6="},
    @{old = "137÷7="; new = "958÷7="},
    @{old = "559÷8="; new = "773÷4="},
    @{old = "643÷6="; new = "423÷8="},
    @{old = "284÷2="; new = "295÷3="},
    @{old = "986÷9="; new = "839÷8="},
    @{old = "981÷4="; new = "113÷3="},
    @{old = "679÷9="; new = "606÷4="},
    @{old = "317÷5="; new = "892÷2="},
    @{old = "381÷2="; new = "319÷8="},
    @{old = "327÷7="; new = "307÷2="},
    @{old = "525÷5="; new = "181÷7="},
    @{old = "393÷9="; new = "184÷6="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
